$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.071804404258728
$ws.Range("B1").Value = 3.077982187271118
$ws.Range("C1").Value = 6.624104499816895
$ws.Range("D1").Value = 1.828423857688904
$ws.Range("E1").Value = 1.275410771369934
